$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header title text (shared strings with multiple runs) ---
$ws.Range("A8").Value = "Volume 32   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/10/2025  Through  2/16/2025"

# --- Donor cells used to carry correct style (number format) when a cell's
#     underlying type changes between numeric and shared-text placeholder ---
$donorText    = $ws.Cells.Item(14, 3)   # C14:  s=13 (General/text style)
$donorNumber  = $ws.Cells.Item(14, 9)   # I14:  s=14 (#,##0 integer style)
$donorPercent = $ws.Cells.Item(15, 13)  # M15:  s=15 (#,##0.0 percent-like style)

# --- Cells that change from a numeric value to a text placeholder ('0' or '***.*') ---
$donorText.Copy($ws.Cells.Item(14, 6))
$ws.Cells.Item(14, 6).Value = "0"
$donorText.Copy($ws.Cells.Item(15, 3))
$ws.Cells.Item(15, 3).Value = "0"
$donorText.Copy($ws.Cells.Item(16, 3))
$ws.Cells.Item(16, 3).Value = "0"
$donorText.Copy($ws.Cells.Item(16, 4))
$ws.Cells.Item(16, 4).Value = "0"
$donorText.Copy($ws.Cells.Item(16, 5))
$ws.Cells.Item(16, 5).Value = "***.*"
$donorText.Copy($ws.Cells.Item(27, 3))
$ws.Cells.Item(27, 3).Value = "0"

# --- Cells that change from a text placeholder to a real numeric value ---
$donorPercent.Copy($ws.Cells.Item(14, 13))
$ws.Cells.Item(14, 13).Value = 0
$donorNumber.Copy($ws.Cells.Item(18, 3))
$ws.Cells.Item(18, 3).Value = 1
$donorPercent.Copy($ws.Cells.Item(29, 13))
$ws.Cells.Item(29, 13).Value = -100
$donorPercent.Copy($ws.Cells.Item(30, 13))
$ws.Cells.Item(30, 13).Value = -100
$donorNumber.Copy($ws.Cells.Item(31, 6))
$ws.Cells.Item(31, 6).Value = 1
$donorNumber.Copy($ws.Cells.Item(31, 9))
$ws.Cells.Item(31, 9).Value = 1

# --- Plain numeric value updates (style unchanged) ---
$ws.Cells.Item(14, 14).Value = -50
$ws.Cells.Item(15, 5).Value = -100
$ws.Cells.Item(15, 7).Value = 2
$ws.Cells.Item(15, 8).Value = -50
$ws.Cells.Item(15, 10).Value = 3
$ws.Cells.Item(15, 11).Value = -33.333333333333
$ws.Cells.Item(15, 14).Value = -60
$ws.Cells.Item(16, 6).Value = 5
$ws.Cells.Item(16, 7).Value = 3
$ws.Cells.Item(16, 8).Value = 66.666666666666
$ws.Cells.Item(16, 12).Value = -9.090909090909
$ws.Cells.Item(16, 13).Value = -75
$ws.Cells.Item(16, 14).Value = -92.753623188405
$ws.Cells.Item(17, 6).Value = 14
$ws.Cells.Item(17, 7).Value = 30
$ws.Cells.Item(17, 8).Value = -53.333333333333
$ws.Cells.Item(17, 9).Value = 25
$ws.Cells.Item(17, 10).Value = 46
$ws.Cells.Item(17, 11).Value = -45.652173913043
$ws.Cells.Item(17, 12).Value = -10.714285714285
$ws.Cells.Item(17, 13).Value = -26.470588235294
$ws.Cells.Item(17, 14).Value = -44.444444444444
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 6
$ws.Cells.Item(18, 7).Value = 3
$ws.Cells.Item(18, 8).Value = 100
$ws.Cells.Item(18, 9).Value = 20
$ws.Cells.Item(18, 10).Value = 16
$ws.Cells.Item(18, 11).Value = 25
$ws.Cells.Item(18, 12).Value = 66.666666666666
$ws.Cells.Item(18, 13).Value = -62.264150943396
$ws.Cells.Item(18, 14).Value = -87.951807228915
$ws.Cells.Item(19, 3).Value = 6
$ws.Cells.Item(19, 5).Value = 20
$ws.Cells.Item(19, 7).Value = 27
$ws.Cells.Item(19, 8).Value = 7.407407407407
$ws.Cells.Item(19, 9).Value = 48
$ws.Cells.Item(19, 10).Value = 50
$ws.Cells.Item(19, 11).Value = -4
$ws.Cells.Item(19, 12).Value = -2.040816326530
$ws.Cells.Item(19, 13).Value = -9.433962264150
$ws.Cells.Item(19, 14).Value = -23.809523809523
$ws.Cells.Item(20, 3).Value = 4
$ws.Cells.Item(20, 4).Value = 5
$ws.Cells.Item(20, 5).Value = -20
$ws.Cells.Item(20, 6).Value = 12
$ws.Cells.Item(20, 7).Value = 14
$ws.Cells.Item(20, 8).Value = -14.285714285714
$ws.Cells.Item(20, 9).Value = 21
$ws.Cells.Item(20, 10).Value = 31
$ws.Cells.Item(20, 11).Value = -32.258064516129
$ws.Cells.Item(20, 12).Value = 50
$ws.Cells.Item(20, 13).Value = -61.111111111111
$ws.Cells.Item(20, 14).Value = -95
$ws.Cells.Item(21, 3).Value = 12
$ws.Cells.Item(21, 4).Value = 20
$ws.Cells.Item(21, 5).Value = -40
$ws.Cells.Item(21, 6).Value = 67
$ws.Cells.Item(21, 7).Value = 79
$ws.Cells.Item(21, 8).Value = -15.189873417721
$ws.Cells.Item(21, 9).Value = 127
$ws.Cells.Item(21, 10).Value = 154
$ws.Cells.Item(21, 11).Value = -17.532467532467
$ws.Cells.Item(21, 12).Value = 9.482758620689
$ws.Cells.Item(21, 13).Value = -46.413502109704
$ws.Cells.Item(21, 14).Value = -84.862932061978
$ws.Cells.Item(24, 3).Value = 18
$ws.Cells.Item(24, 4).Value = 18
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = 54
$ws.Cells.Item(24, 7).Value = 71
$ws.Cells.Item(24, 8).Value = -23.943661971831
$ws.Cells.Item(24, 9).Value = 86
$ws.Cells.Item(24, 10).Value = 103
$ws.Cells.Item(24, 11).Value = -16.504854368932
$ws.Cells.Item(24, 12).Value = -23.214285714285
$ws.Cells.Item(24, 13).Value = -15.686274509803
$ws.Cells.Item(25, 3).Value = 5
$ws.Cells.Item(25, 4).Value = 3
$ws.Cells.Item(25, 5).Value = 66.666666666666
$ws.Cells.Item(25, 7).Value = 16
$ws.Cells.Item(25, 8).Value = -31.25
$ws.Cells.Item(25, 9).Value = 20
$ws.Cells.Item(25, 10).Value = 28
$ws.Cells.Item(25, 11).Value = -28.571428571428
$ws.Cells.Item(25, 12).Value = -16.666666666666
$ws.Cells.Item(26, 3).Value = 11
$ws.Cells.Item(26, 4).Value = 11
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 36
$ws.Cells.Item(26, 7).Value = 22
$ws.Cells.Item(26, 8).Value = 63.636363636363
$ws.Cells.Item(26, 9).Value = 61
$ws.Cells.Item(26, 10).Value = 40
$ws.Cells.Item(26, 11).Value = 52.5
$ws.Cells.Item(26, 12).Value = 74.285714285714
$ws.Cells.Item(26, 13).Value = -18.666666666666
$ws.Cells.Item(27, 5).Value = -100
$ws.Cells.Item(27, 10).Value = 5
$ws.Cells.Item(27, 11).Value = -60
$ws.Cells.Item(27, 12).Value = -33.333333333333
$ws.Cells.Item(28, 3).Value = 2
$ws.Cells.Item(28, 4).Value = 3
$ws.Cells.Item(28, 5).Value = -33.333333333333
$ws.Cells.Item(28, 6).Value = 3
$ws.Cells.Item(28, 7).Value = 6
$ws.Cells.Item(28, 8).Value = -50
$ws.Cells.Item(28, 9).Value = 8
$ws.Cells.Item(28, 10).Value = 8
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 166.666666666667
